# ZBP_09_jednotlive_aktivity.xlsx update: add new survey wave "30. 3. 2021"
# and refresh the "aktualizace" (last updated) date in the two title cells.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "data" - percentages, new column AB for the 30. 3. 2021 wave
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("data")

$ws1.Range("AB1").Value = "30. 3. 2021"
$ws1.Range("AA1").Copy()
$ws1.Range("AB1").PasteSpecial(-4122)

$ws1.Range("AB2").Value = 0.87
$ws1.Range("AB3").Value = 0.68
$ws1.Range("AB4").Value = 0.63
$ws1.Range("AB5").Value = 0.64
$ws1.Range("AB6").Value = 0.61
$ws1.Range("AB7").Value = 0.35
$ws1.Range("AB8").Value = 0.83
$ws1.Range("AB9").Value = 0.84
$ws1.Range("AB10").Value = 0.92
$ws1.Range("AB11").Value = 0.92
$ws1.Range("AB12").Value = 0.91
$ws1.Range("AB13").Value = 0.83
$ws1.Range("AB14").Value = 0.72
$ws1.Range("AB15").Value = 0.88
$ws1.Range("AB16").Value = 0.93
$ws1.Range("AB17").Value = 0.85
$ws1.Range("AB18").Value = 0.93
$ws1.Range("AB19").Value = 0.83
$ws1.Range("AB20").Value = 0.85
$ws1.Range("AB21").Value = 0.89
$ws1.Range("AB22").Value = 0.87
$ws1.Range("AB23").Value = 0.88
$ws1.Range("AB24").Value = 0.86
$ws1.Range("AB25").Value = 0.88
$ws1.Range("AB26").Value = 0.57
$ws1.Range("AB27").Value = 0.66
$ws1.Range("AB28").Value = 0.77
$ws1.Range("AB29").Value = 0.75
$ws1.Range("AB30").Value = 0.72
$ws1.Range("AB31").Value = 0.61
$ws1.Range("AB32").Value = 0.48
$ws1.Range("AB33").Value = 0.67
$ws1.Range("AB34").Value = 0.78
$ws1.Range("AB35").Value = 0.67
$ws1.Range("AB36").Value = 0.75
$ws1.Range("AB37").Value = 0.62
$ws1.Range("AB38").Value = 0.65
$ws1.Range("AB39").Value = 0.72
$ws1.Range("AB40").Value = 0.6899999999999999
$ws1.Range("AB41").Value = 0.7
$ws1.Range("AB42").Value = 0.61
$ws1.Range("AB43").Value = 0.68
$ws1.Range("AB44").Value = 0.48
$ws1.Range("AB45").Value = 0.61
$ws1.Range("AB46").Value = 0.74
$ws1.Range("AB47").Value = 0.74
$ws1.Range("AB48").Value = 0.66
$ws1.Range("AB49").Value = 0.54
$ws1.Range("AB50").Value = 0.42
$ws1.Range("AB51").Value = 0.6
$ws1.Range("AB52").Value = 0.75
$ws1.Range("AB53").Value = 0.6
$ws1.Range("AB54").Value = 0.73
$ws1.Range("AB55").Value = 0.5600000000000001
$ws1.Range("AB56").Value = 0.61
$ws1.Range("AB57").Value = 0.66
$ws1.Range("AB58").Value = 0.62
$ws1.Range("AB59").Value = 0.66
$ws1.Range("AB60").Value = 0.64
$ws1.Range("AB61").Value = 0.64
$ws1.Range("AB62").Value = 0.49
$ws1.Range("AB63").Value = 0.62
$ws1.Range("AB64").Value = 0.75
$ws1.Range("AB65").Value = 0.75
$ws1.Range("AB66").Value = 0.6899999999999999
$ws1.Range("AB67").Value = 0.51
$ws1.Range("AB68").Value = 0.47
$ws1.Range("AB69").Value = 0.61
$ws1.Range("AB70").Value = 0.74
$ws1.Range("AB71").Value = 0.6
$ws1.Range("AB72").Value = 0.74
$ws1.Range("AB73").Value = 0.58
$ws1.Range("AB74").Value = 0.61
$ws1.Range("AB75").Value = 0.67
$ws1.Range("AB76").Value = 0.64
$ws1.Range("AB77").Value = 0.67
$ws1.Range("AB78").Value = 0.62
$ws1.Range("AB79").Value = 0.63
$ws1.Range("AB80").Value = 0.48
$ws1.Range("AB81").Value = 0.57
$ws1.Range("AB82").Value = 0.72
$ws1.Range("AB83").Value = 0.6899999999999999
$ws1.Range("AB84").Value = 0.65
$ws1.Range("AB85").Value = 0.54
$ws1.Range("AB86").Value = 0.43
$ws1.Range("AB87").Value = 0.58
$ws1.Range("AB88").Value = 0.7
$ws1.Range("AB89").Value = 0.58
$ws1.Range("AB90").Value = 0.7
$ws1.Range("AB91").Value = 0.54
$ws1.Range("AB92").Value = 0.58
$ws1.Range("AB93").Value = 0.63
$ws1.Range("AB94").Value = 0.58
$ws1.Range("AB95").Value = 0.63
$ws1.Range("AB96").Value = 0.64
$ws1.Range("AB97").Value = 0.66
$ws1.Range("AB98").Value = 0.24
$ws1.Range("AB99").Value = 0.27
$ws1.Range("AB100").Value = 0.49
$ws1.Range("AB101").Value = 0.52
$ws1.Range("AB102").Value = 0.33
$ws1.Range("AB103").Value = 0.2
$ws1.Range("AB104").Value = 0.24
$ws1.Range("AB105").Value = 0.29
$ws1.Range("AB106").Value = 0.44
$ws1.Range("AB107").Value = 0.3
$ws1.Range("AB108").Value = 0.39
$ws1.Range("AB109").Value = 0.34
$ws1.Range("AB110").Value = 0.35
$ws1.Range("AB111").Value = 0.35
$ws1.Range("AB112").Value = 0.38
$ws1.Range("AB113").Value = 0.33
$ws1.Range("AB114").Value = 0.34
$ws1.Range("AB115").Value = 0.26

$ws1.Range("A116").Value = "Život během pandemie, Jednotlivé protektivní aktivity, % respondentů celkově a ve skupinách, aktualizace 7. 4. 2021"

# ---------------------------------------------------------------------
# Sheet 2: "pocetR" - respondent counts, new column AA for the wave
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pocetR")

$ws2.Range("AA1").Value = "30. 3. 2021"
$ws2.Range("Z1").Copy()
$ws2.Range("AA1").PasteSpecial(-4122)

$ws2.Range("AA2").Value = 2061
$ws2.Range("AA3").Value = 488
$ws2.Range("AA4").Value = 758
$ws2.Range("AA5").Value = 815
$ws2.Range("AA6").Value = 627
$ws2.Range("AA7").Value = 715
$ws2.Range("AA8").Value = 498
$ws2.Range("AA9").Value = 387
$ws2.Range("AA10").Value = 751
$ws2.Range("AA11").Value = 923
$ws2.Range("AA12").Value = 626
$ws2.Range("AA13").Value = 734
$ws2.Range("AA14").Value = 701
$ws2.Range("AA15").Value = 1007
$ws2.Range("AA16").Value = 1054
$ws2.Range("AA17").Value = 1084
$ws2.Range("AA18").Value = 464
$ws2.Range("AA19").Value = 241
$ws2.Range("AA20").Value = 272

# Existing value revised in this update (21. 7. 2020 wave, row 8)
$ws2.Range("Z8").Value = 435

$ws2.Range("AA21").Value = ""

$ws2.Range("A21").Value = "Život během pandemie, Jednotlivé protektivní aktivity, velikost dotázaného souboru celkově a ve skupinách, aktualizace 7. 4. 2021"
